$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '37.492.09'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '2.065.42'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Formula = "'231.84"
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Formula = "'57.62"
$ws.Range("E8").Value = '  -2.50%  '
$ws.Range("D9").Formula = "'0.387"
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("E10").Value = '  -1.46%  '
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("D12").Value = '2.368.89'
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").Formula = "'14.71"
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").Formula = "'21.33"
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("D16").Formula = "'5.36"
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").Value = '2.061.11'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").Value = '37.493.02'
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("D19").Formula = "'6.16"
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Formula = "'69.83"
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("E21").Value = '  -2.43%  '
$ws.Range("D22").Formula = "'227.27"
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E25").Value = '  -2.25%  '
$ws.Range("D26").Formula = "'9.91"
$ws.Range("E26").Value = '  +7.56%  '
$ws.Range("D27").Formula = "'169.90"
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("E28").Value = '  -4.47%  '
$ws.Range("D29").Formula = "'19.25"
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("E30").Value = '  -4.13%  '
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("D32").Formula = "'4.56"
$ws.Range("E32").Value = '  -3.51%  '
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  -3.37%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("E40").Value = '  +4.07%  '
$ws.Range("D41").Formula = "'98.78"
$ws.Range("E41").Value = '  -0.33%  '
$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").Formula = "'0.0962"
$ws.Range("E42").Value = '  -2.04%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Formula = "'1.20"
$ws.Range("E43").Value = '  +4.51%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '1.477.03'
$ws.Range("E45").Value = '  +2.15%  '
$ws.Range("D46").Formula = "'16.73"
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("E48").Value = '  -4.49%  '
$ws.Range("E49").Value = '  -1.77%  '
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").Value = '2.253.11'
$ws.Range("E51").Value = '  -1.04%  '
